$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (save games filtered), columns B:E plus
# recomputed G (sum of B:E). Column F (Win) is left untouched.

$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    4 = @{ B = 1.459612070389937;  C = 1.667794583268128; D = 0.1575252929769615; E = 8.660232485948974; G = 11.945164432584 }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 26.21740644021617;  E = 0.496779210170732; G = 31.61296591696135 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
